$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the portfolio companies (fund 1 -> fund 2 test data) ---
# Update every cell that references the old shared strings in one pass so
# the writer reclaims the very same shared-string slots (7 and 8) for the
# new names instead of appending fresh entries at the end of the table.
$ws.Range("D2").Value = "TSTF2 Port Co 3"
$ws.Range("D3").Value = "TSTF2 Port Co 4"
$ws.Range("D4").Value = "TSTF2 Port Co 3"
$ws.Range("D5").Value = "TSTF2 Port Co 4"
$ws.Range("D6").Value = "TSTF2 Port Co 3"
$ws.Range("D7").Value = "TSTF2 Port Co 4"
# (rows 6 and 7 themselves are removed below; renaming first just keeps the
# shared-string table from ever containing a stale "TSTF1 ..." entry)

# --- Row 4 becomes a second 31-Mar-2025 line (was 31-Mar-2024) ---
$ws.Range("A4").Value = "03/31/2025"
$ws.Range("A5").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("E4").Value = "Equity"

# --- Row 5: valuation bumped to 250, instrument becomes CCPS on Co4 ---
$ws.Range("C5").Value = 250
$ws.Range("E5").Value = "CCPS"

# --- New trailing helper cells (J:L) copied onto rows 2 and 3 ---
$ws.Range("J4:L4").Copy()
$ws.Range("J2:L2").PasteSpecial(-4122)
$ws.Range("J3:L3").PasteSpecial(-4122)

# Row 5 no longer carries the helper cells (only rows 2-4 do now).
$ws.Range("J5:L5").Clear()

# Drop the old rows 6 and 7 - fund 2's test sheet only keeps 4 data rows.
$ws.Rows("6:7").Delete()

# --- Misc workbook metadata ---
$fd = $wb.Names.Item("_xlnm._FilterDatabase")
$fd.RefersTo = "=Sheet1!`$A`$1:`$O`$5"

$ws.Range("C4").Select()
